$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $rng = $d.Content
    $n = 0
    while ($rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false)) {
        $rng.Text = $replace
        $n += 1
        $rng.Collapse(0)
    }
    return $n
}

# 1. Title change (occurs twice: Heading1 at the top, and bold run near the bottom)
Replace-All "Play Break Bones for Free - Review of Hacksaw Gaming's Slot Game" "Play Break Bones Free - Review of Hacksaw Gaming's Video Slot"

# 2. "What we like" bullet points
Replace-All "Appealing Urban Graffiti theme" "Eye-catching graphics and Urban Graffiti theme"
Replace-All "Medium volatility and high RTP (96.22%)" "Medium volatility and high RTP"
Replace-All "Advanced settings for automatic spins" "Advanced settings for betting and automatic spins"
Replace-All "Global Multiplier feature for bigger winnings" "Multiple special features including Wild multipliers and Echo Spins"

# 3. "What we don't like" bullet point
Replace-All "No progressive jackpot feature" "Standard paytable symbols may feel repetitive after extended gameplay"

# 4. Meta description paragraph (italic run near the bottom)
Replace-All "Read our review of Break Bones, a 3-reel, 17-fixed-payline video slot game by Hacksaw Gaming. Play for free and learn about its special features." "Read our review of Break Bones, an Urban Graffiti-themed video slot game. Play for free and enjoy multiple special features."
